$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date in A1 by one day
$ws.Range("A1:E1").UnMerge()
$ws.Range("A1").Value = 45309
$ws.Range("A1:E1").Merge()

# Step 2: update pricing table values
$ws.Range("D23").Value = 11284.427
$ws.Range("D24").Value = 15555.424
$ws.Range("D25").Value = 22119.28
$ws.Range("D26").Value = 32999.091
$ws.Range("D27").Value = 59524.242
$ws.Range("D28").Value = 77057.803
$ws.Range("D36").Value = 6250
$ws.Range("D37").Value = 9150
